# Memphis roster housekeeping edit:
# Rows 7 and 8 (Ja Morant / David Roddy) swap places in the roster table -
# every column except the leading row index (column A) trades places
# between the two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 7 becomes the old "David Roddy" row ---
$ws.Range("B7").Value = 27
$ws.Range("C7").Value = "David Roddy"
$ws.Range("D7").Value = "PF"
$ws.Range("E7").Value = "6-6"
$ws.Range("F7").Value = 255
$ws.Range("G7").Value = "March 27, 2001"
$ws.Range("H7").Value = "us"
$ws.Range("I7").Value = "R"
$ws.Range("J7").Value = "Colorado State"
$ws.Range("K7").Value = "https://www.basketball-reference.com/players/r/roddyda01.html"

# --- Row 8 becomes the old "Ja Morant" row ---
$ws.Range("B8").Value = 12
$ws.Range("C8").Value = "Ja Morant"
$ws.Range("D8").Value = "PG"
$ws.Range("E8").Value = "6-3"
$ws.Range("F8").Value = 174
$ws.Range("G8").Value = "August 10, 1999"
$ws.Range("H8").Value = "us"
$ws.Range("I8").Value = "'3"
$ws.Range("J8").Value = "Murray State"
$ws.Range("K8").Value = "https://www.basketball-reference.com/players/m/moranja01.html"
